# Updates cryptos list (Price / Volume(1h) columns) to match the latest snapshot.
# Values are written as literal text (matching the sheet's existing inline-string
# cells): numeric-looking Price values are forced to Text via NumberFormat "@" so
# Excel does not auto-convert them to real numbers, then the style is reset back to
# "Normal" so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.213.52"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "2.521.00"
$ws.Range("E3").Value = "  +3.06%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.573"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("D9").Value = "2.543.27"
$ws.Range("E9").Value = "  +3.27%  "

$ws.Range("E10").Value = "  +1.42%  "

$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.58%  "

$ws.Range("E13").Value = "  +2.63%  "

$ws.Range("D14").Value = "2.964.85"
$ws.Range("E14").Value = "  +2.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.90%  "

$ws.Range("D16").Value = "59.143.58"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("E17").Value = "  +2.21%  "

$ws.Range("D18").Value = "2.537.61"
$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.34%  "

$ws.Range("E20").Value = "  -1.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.42%  "

$ws.Range("E23").Value = "  +2.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.433"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.94%  "

$ws.Range("E26").Value = "  +3.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.03%  "

$ws.Range("E28").Value = "  +3.90%  "

$ws.Range("D29").Value = "0.0₃0780"
$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.04%  "

$ws.Range("E33").Value = "  +7.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.20%  "

$ws.Range("E38").Value = "  -5.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "296.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.826"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.602"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.63%  "

$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0516"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.06%  "

$ws.Range("E51").Value = "  -0.69%  "

